$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for specific rows as per repull/recalculation
$ws.Range("F10").Value = -1
$ws.Range("F12").Value = -6
$ws.Range("F13").Value = -5
$ws.Range("F16").Value = -7
$ws.Range("F17").Value = -3
$ws.Range("F19").Value = -5
$ws.Range("F23").Value = -6
$ws.Range("F25").Value = 0
